$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 11 (Fats Waller) values: C11, D11, E11
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 21
$ws.Range("E11").Value = 1904

# Delete rows 15 and 16 (the two trailing test rows) entirely, shifting cells up
$ws.Range("A15:E16").Delete()

# Update the selection to match the target workbook state
$ws.Range("A16").Select()
